$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.971.40'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.19%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.549.71'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.60%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '547.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.96%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.34%  '

$ws.Range('E7').Value = '  +0.38%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.81%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.546.85'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.97%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.67'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.19%  '

$ws.Range('E11').Value = '  -1.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.161'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.21%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.350'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.11%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.015.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.91%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.020.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.91%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.75%  '

$ws.Range('E17').Value = '  -1.62%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.569.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.41%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '334.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.79%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.66%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.36%  '

$ws.Range('E23').Value = '  -0.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.473'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.39%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '62.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.73%  '

$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('E27').Value = '  -3.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.14%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0756'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.92%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.26%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.82%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.78%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.94'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.87%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.05'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.09%  '

$ws.Range('E36').Value = '  +0.87%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.882'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.85%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.30'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.838'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.80%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.64'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.27%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '283.53'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.64%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '135.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.58%  '

$ws.Range('E44').Value = '  +0.86%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0965'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.44%  '

$ws.Range('E46').Value = '  -0.17%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.584'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.47%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0525'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.79%  '

$ws.Range('E49').Value = '  -1.90%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.944.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.84%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.26%  '
